$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated roster data (Oyuncu Adi, Pozisyon, Takim)
$data = @(
    @("Cade Cunningham", "PG,SG", "Detroit Pistons"),
    @("Carlton Carrington", "PG,SG", "Washington Wizards"),
    @("Derrick White", "PG,SG", "Boston Celtics"),
    @("Malik Monk", "PG,SG,SF", "Sacramento Kings"),
    @("Ausar Thompson", "SF,PF", "Detroit Pistons"),
    @("Precious Achiuwa", "PF,C", "New York Knicks"),
    @("Santi Aldama", "PF,C", "Memphis Grizzlies"),
    @("Naz Reid", "PF,C", "Minnesota Timberwolves"),
    @("Isaiah Hartenstein", "C", "Oklahoma City Thunder"),
    @("Jaxson Hayes", "PF,C", "Los Angeles Lakers"),
    @("Onyeka Okongwu", "PF,C", "Atlanta Hawks"),
    @("Damian Lillard", "PG", "Milwaukee Bucks"),
    @("Coby White", "PG,SG", "Chicago Bulls"),
    @("Devin Vassell", "SG,SF", "San Antonio Spurs"),
    @("Andrew Wiggins", "SF,PF", "Miami Heat"),
    @("Collin Sexton", "PG,SG", "Utah Jazz"),
    @("LaMelo Ball", "PG,SG", "Charlotte Hornets"),
    @("Anthony Davis", "PF,C", "Dallas Mavericks")
)

$row = 2
foreach ($player in $data) {
    $ws.Cells.Item($row, 1).Value = $player[0]
    $ws.Cells.Item($row, 2).Value = $player[1]
    $ws.Cells.Item($row, 3).Value = $player[2]
    $row = $row + 1
}
